# daily auto push: 2026-01-12 18:48 UTC
# Insert two new daily-log rows (2026/01/12 22:00 and 2026/01/13 01:00)
# right after the existing 2026/01/12 entries, pushing the 2026/12/29 ...
# 2027/01/05 block that follows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything currently at row 637 onward down by two rows.
$ws.Rows("637:638").Insert()

# Row 637: 2026/01/12, 月, 22, 187
$ws.Range("A637").Value = "'2026/01/12"
$ws.Range("A637").ClearFormats()
$ws.Range("B637").Value = "月"
$ws.Range("C637").Value = 22
$ws.Range("D637").Value = 187

# Row 638: 2026/01/13, 火, 1, 201
$ws.Range("A638").Value = "'2026/01/13"
$ws.Range("A638").ClearFormats()
$ws.Range("B638").Value = "火"
$ws.Range("C638").Value = 1
$ws.Range("D638").Value = 201
